# Auto-generated edit script applying cryptos.xlsx diff (97 cell changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.918.39"
$ws.Range("E2").Value = "  -4.84%  "
$ws.Range("D3").Value = "3.540.89"
$ws.Range("E3").Value = "  -5.67%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "3.530.25"
$ws.Range("E7").Value = "  -5.80%  "
$ws.Range("E8").Value = "  -5.50%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.660"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.82%  "
$ws.Range("E11").Value = "  -12.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -15.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.52%  "
$ws.Range("D15").Value = "4.101.92"
$ws.Range("E15").Value = "  -5.84%  "
$ws.Range("D16").Value = "3.537.07"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.02%  "
$ws.Range("D20").Value = "65.633.50"
$ws.Range("E20").Value = "  -5.06%  "
$ws.Range("E21").Value = "  -8.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.11%  "
$ws.Range("E25").Value = "  -7.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "611.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "62.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.90%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.05%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -8.86%  "
$ws.Range("E40").Value = "  -18.88%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.32%  "
$ws.Range("D43").Value = "2.873.27"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0404"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.67%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.41"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.23%  "
